# 9 Mayis 2020 verileri eklendi
# Appends the 9 May 2020 COVID-19 data row (row 59) to the "data" sheet,
# extends Table3 to cover the new row, and updates the sheet's
# selection/dimension to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: date, test, case, death, recovered
$ws.Range("A59").Value = 43960
$ws.Range("B59").Value = 35605
$ws.Range("B59").NumberFormat = "General"
$ws.Range("C59").Value = 1546
$ws.Range("D59").Value = 50
$ws.Range("E59").Value = 3084

# Grow the table (and its autofilter) to include the newly added row.
$ws.ListObjects("Table3").Resize($ws.Range("A1:E59"))

# Match the author's last selected cell after entering the new data.
$ws.Range("E58").Select()
